$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.654881412372333
$ws.Range("C2").Value = 0.3269421160264017
$ws.Range("D2").Value = 0.04116341890107122
$ws.Range("E2").Value = 0.0709066865514405
$ws.Range("F2").Value = 5.035716064856842
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1825092295516555
$ws.Range("M2").Value = 0.4906818933942318
$ws.Range("B3").Value = 1.58751167872316
$ws.Range("C3").Value = 0.3113641504563986
$ws.Range("D3").Value = 0.03808527330661349
$ws.Range("E3").Value = 0.07151317091575393
$ws.Range("F3").Value = 4.851349761972102
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1800861087912082
$ws.Range("M3").Value = 0.4803755351055017
$ws.Range("B4").Value = 1.547709735660533
$ws.Range("C4").Value = 0.3021723861938312
$ws.Range("D4").Value = 0.03620435605168382
$ws.Range("E4").Value = 0.0719124968995537
$ws.Range("F4").Value = 4.739545318432931
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1786455481355631
$ws.Range("M4").Value = 0.4745060722007537
$ws.Range("B5").Value = 1.531880906409071
$ws.Range("C5").Value = 0.29851977324671
$ws.Range("D5").Value = 0.03543994555714391
$ws.Range("E5").Value = 0.07208201892930788
$ws.Range("F5").Value = 4.694327541323332
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1780702718664458
$ws.Range("M5").Value = 0.4722291538982688
$ws.Range("B6").Value = 1.529276077001498
$ws.Range("C6").Value = 0.2979188619858917
$ws.Range("D6").Value = 0.03531313537812508
$ws.Range("E6").Value = 0.07211057879040172
$ws.Range("F6").Value = 4.686839677558993
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1779754548215138
$ws.Range("M6").Value = 0.4718580036487126
$ws.Range("B7").Value = 1.547494683280604
$ws.Range("C7").Value = 0.3021227497483778
$ws.Range("D7").Value = 0.0361940387962818
$ws.Range("E7").Value = 0.07191475560156402
$ws.Range("F7").Value = 4.738934115475018
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1786377422792498
$ws.Range("M7").Value = 0.4744749000346147
$ws.Range("B8").Value = 1.631326407921506
$ws.Range("C8").Value = 0.3214929308545038
$ws.Range("D8").Value = 0.04010006669503241
$ws.Range("E8").Value = 0.07111022167729963
$ws.Range("F8").Value = 4.97185173938027
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1816638604472516
$ws.Range("M8").Value = 0.4870328087104596
$ws.Range("B9").Value = 1.808243947042683
$ws.Range("C9").Value = 0.362475079288032
$ws.Range("D9").Value = 0.0478413228036203
$ws.Range("E9").Value = 0.0697454833780915
$ws.Range("F9").Value = 5.440070794699466
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.187978610368674
$ws.Range("M9").Value = 0.5153194574813611
$ws.Range("B10").Value = 1.946046021915095
$ws.Range("C10").Value = 0.3944671454910633
$ws.Range("D10").Value = 0.05359221827968952
$ws.Range("E10").Value = 0.06887151593481544
$ws.Range("F10").Value = 5.791629124467505
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1928585155893998
$ws.Range("M10").Value = 0.5383669275127048
$ws.Range("B11").Value = 2.010476537507486
$ws.Range("C11").Value = 0.4094424976647701
$ws.Range("D11").Value = 0.05622521823730153
$ws.Range("E11").Value = 0.06850163978817303
$ws.Range("F11").Value = 5.953329208850789
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1951326180401836
$ws.Range("M11").Value = 0.5493515098151605
$ws.Range("B12").Value = 2.035128669714084
$ws.Range("C12").Value = 0.4151749363073804
$ws.Range("D12").Value = 0.05722493247122884
$ws.Range("E12").Value = 0.06836554187139932
$ws.Range("F12").Value = 6.014825675715656
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1960017025121346
$ws.Range("M12").Value = 0.5535836053970229
$ws.Range("B13").Value = 2.029808060759478
$ws.Range("C13").Value = 0.4139375986554512
$ws.Range("D13").Value = 0.05700950445249475
$ws.Range("E13").Value = 0.0683946768575181
$ws.Range("F13").Value = 6.001569418181191
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1958141746858431
$ws.Range("M13").Value = 0.5526689156965787
$ws.Range("B14").Value = 2.012499581583029
$ws.Range("C14").Value = 0.4099128689960878
$ws.Range("D14").Value = 0.05630741083128044
$ws.Range("E14").Value = 0.06849036353676041
$ws.Range("F14").Value = 5.958383214121511
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1952039584008176
$ws.Range("M14").Value = 0.5496982311914635
$ws.Range("B15").Value = 2.001930766140333
$ws.Range("C15").Value = 0.4074556553162267
$ws.Range("D15").Value = 0.05587771062630509
$ws.Range("E15").Value = 0.06854949041490421
$ws.Range("F15").Value = 5.931965087199615
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1948312202674813
$ws.Range("M15").Value = 0.5478880558615131
$ws.Range("B16").Value = 1.941870673897483
$ws.Range("C16").Value = 0.3934970403219609
$ws.Range("D16").Value = 0.05342050462101611
$ws.Range("E16").Value = 0.06889624410961837
$ws.Range("F16").Value = 5.78109811528833
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1927110000424932
$ws.Range("M16").Value = 0.5376591703948108
$ws.Range("B17").Value = 1.905474422068608
$ws.Range("C17").Value = 0.3850426090877193
$ws.Range("D17").Value = 0.05191757373354733
$ws.Range("E17").Value = 0.06911604807206029
$ws.Range("F17").Value = 5.689006943400841
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1914242938555404
$ws.Range("M17").Value = 0.5315125704641659
$ws.Range("B18").Value = 1.884704199872147
$ws.Range("C18").Value = 0.3802195082721482
$ws.Range("D18").Value = 0.05105470363834286
$ws.Range("E18").Value = 0.06924508135166185
$ws.Range("F18").Value = 5.636205009110824
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1906893086200014
$ws.Range("M18").Value = 0.5280242533332711
$ws.Range("B19").Value = 1.877699827562878
$ws.Range("C19").Value = 0.3785932741721751
$ws.Range("D19").Value = 0.05076281449517239
$ws.Range("E19").Value = 0.06928921821064815
$ws.Range("F19").Value = 5.618355501492175
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1904413260675568
$ws.Range("M19").Value = 0.5268512326574424
$ws.Range("B20").Value = 1.909331876987892
$ws.Range("C20").Value = 0.3859384860682553
$ws.Range("D20").Value = 0.05207739855495674
$ws.Range("E20").Value = 0.06909237978237748
$ws.Range("F20").Value = 5.698792911956531
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.191560737724771
$ws.Range("M20").Value = 0.5321620136780325
$ws.Range("B21").Value = 2.017576592263765
$ws.Range("C21").Value = 0.4110933504912282
$ws.Range("D21").Value = 0.0565135588982173
$ws.Range("E21").Value = 0.0684621505296823
$ws.Range("F21").Value = 5.97106080764371
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1953829773679558
$ws.Range("M21").Value = 0.5505688218808302
$ws.Range("B22").Value = 2.089800933539891
$ws.Range("C22").Value = 0.4278929910621514
$ws.Range("D22").Value = 0.05942842401701398
$ws.Range("E22").Value = 0.06807336955785281
$ws.Range("F22").Value = 6.150547587691165
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1979273261513654
$ws.Range("M22").Value = 0.5630214304542775
$ws.Range("B23").Value = 2.051117007217215
$ws.Range("C23").Value = 0.4188935049157863
$ws.Range("D23").Value = 0.05787120500758647
$ws.Range("E23").Value = 0.06827876008001788
$ws.Range("F23").Value = 6.054607832565921
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1965650777901473
$ws.Range("M23").Value = 0.5563363756703197
$ws.Range("B24").Value = 1.907587440708539
$ws.Range("C24").Value = 0.3855333434611623
$ws.Range("D24").Value = 0.05200513807970708
$ws.Range("E24").Value = 0.06910307191386433
$ws.Range("F24").Value = 5.69436823237271
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1914990366374525
$ws.Range("M24").Value = 0.531868258891059
$ws.Range("B25").Value = 1.759023393854477
$ws.Range("C25").Value = 0.3510618189814636
$ws.Range("D25").Value = 0.04573697145602296
$ws.Range("E25").Value = 0.07009199948946243
$ws.Range("F25").Value = 5.312119290078215
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1862286961929769
$ws.Range("M25").Value = 0.5072718644532017
